$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header date text
$ws.Range("C1").Value = "Amount of Workshop Items as of 2024-04-22"

$data = New-Object 'object[,]' 100,2
$data[0,0] = "Counter-Strike 2"
$data[0,1] = 4648
$data[1,0] = "Dota 2"
$data[1,1] = 32603
$data[2,0] = "Wallpaper Engine"
$data[2,1] = 2187431
$data[3,0] = "Rust"
$data[3,1] = 118623
$data[4,0] = "Team Fortress 2"
$data[4,1] = 9527
$data[5,0] = "Warframe"
$data[5,1] = 881
$data[6,0] = "Sid Meier’s Civilization® VI"
$data[6,1] = 9316
$data[7,0] = "Unturned"
$data[7,1] = 101670
$data[8,0] = "RimWorld"
$data[8,1] = 35691
$data[9,0] = "Hearts of Iron IV"
$data[9,1] = 48550
$data[10,0] = "Don't Starve Together"
$data[10,1] = 17179
$data[11,0] = "DayZ"
$data[11,1] = 58704
$data[12,0] = "tModLoader"
$data[12,1] = 6960
$data[13,0] = "Euro Truck Simulator 2"
$data[13,1] = 22381
$data[14,0] = "Myth of Empires"
$data[14,1] = 53
$data[15,0] = "ARK: Survival Evolved"
$data[15,1] = 19625
$data[16,0] = "Squad"
$data[16,1] = 463
$data[17,0] = "Mount & Blade II: Bannerlord"
$data[17,1] = 619
$data[18,0] = "Slay the Spire"
$data[18,1] = 947
$data[19,0] = "Left 4 Dead 2"
$data[19,1] = 142550
$data[20,0] = "Project Zomboid"
$data[20,1] = 26808
$data[21,0] = "Garry's Mod"
$data[21,1] = 1807968
$data[22,0] = "Crosshair X"
$data[22,1] = "N/A"
$data[23,0] = "VPet-Simulator"
$data[23,1] = 601
$data[24,0] = "Europa Universalis IV"
$data[24,1] = 12834
$data[25,0] = "Crusader Kings III"
$data[25,1] = 8235
$data[26,0] = "Cities: Skylines"
$data[26,1] = 340432
$data[27,0] = "Rocket League®"
$data[27,1] = 608
$data[28,0] = "VTube Studio"
$data[28,1] = 874
$data[29,0] = "Stellaris"
$data[29,1] = 29413
$data[30,0] = "Cookie Clicker"
$data[30,1] = 1245
$data[31,0] = "Arma 3"
$data[31,1] = 144321
$data[32,0] = "Aimlabs"
$data[32,1] = 42675
$data[33,0] = "Golf With Your Friends"
$data[33,1] = 19108
$data[34,0] = "YoloMouse - Game Cursor Changer"
$data[34,1] = 125
$data[35,0] = "Conan Exiles"
$data[35,1] = 3902
$data[36,0] = "Brotato"
$data[36,1] = 251
$data[37,0] = "Victoria 3"
$data[37,1] = 4861
$data[38,0] = "鬼谷八荒 Tale of Immortal"
$data[38,1] = 9103
$data[39,0] = "Kenshi"
$data[39,1] = 14736
$data[40,0] = "DSX"
$data[40,1] = 2
$data[41,0] = "Mount & Blade: Warband"
$data[41,1] = 401
$data[42,0] = "Noita"
$data[42,1] = 1255
$data[43,0] = "American Truck Simulator"
$data[43,1] = 10727
$data[44,0] = "Dying Light"
$data[44,1] = 806
$data[45,0] = "Company of Heroes 2"
$data[45,1] = 14909
$data[46,0] = "Planet Zoo"
$data[46,1] = 104266
$data[47,0] = "People Playground"
$data[47,1] = 709113
$data[48,0] = "SAO Utils 2: Progressive"
$data[48,1] = "N/A"
$data[49,0] = "Divinity: Original Sin 2 - Definitive Edition"
$data[49,1] = 4704
$data[50,0] = "XCOM® 2"
$data[50,1] = 8543
$data[51,0] = "They Are Billions"
$data[51,1] = 4302
$data[52,0] = "Space Engineers"
$data[52,1] = 552484
$data[53,0] = "Halo: The Master Chief Collection"
$data[53,1] = 1143
$data[54,0] = "Dead Cells"
$data[54,1] = 724
$data[55,0] = "F1® 23"
$data[55,1] = 1284
$data[56,0] = "Age of Empires II (Retired)"
$data[56,1] = 17415
$data[57,0] = "Tabletop Simulator"
$data[57,1] = 82664
$data[58,0] = "Hero's Adventure: Road to Passion"
$data[58,1] = 771
$data[59,0] = "Transport Fever 2"
$data[59,1] = 13049
$data[60,0] = "Call of Duty®: Black Ops III"
$data[60,1] = 5255
$data[61,0] = "Farthest Frontier"
$data[61,1] = "N/A"
$data[62,0] = "Warhammer 40,000: Rogue Trader"
$data[62,1] = 14
$data[63,0] = "Killing Floor 2"
$data[63,1] = 3012
$data[64,0] = "A Dance of Fire and Ice"
$data[64,1] = 18870
$data[65,0] = "MyDockFinder"
$data[65,1] = 3660
$data[66,0] = "X4: Foundations"
$data[66,1] = 829
$data[67,0] = "Total War: WARHAMMER II"
$data[67,1] = 12719
$data[68,0] = "SAO Utils: Beta"
$data[68,1] = 263
$data[69,0] = "Pummel Party"
$data[69,1] = 1353
$data[70,0] = "觅长生"
$data[70,1] = 1512
$data[71,0] = "Company of Heroes 3"
$data[71,1] = 509
$data[72,0] = "Human Fall Flat"
$data[72,1] = 504748
$data[73,0] = "Call to Arms - Gates of Hell: Ostfront"
$data[73,1] = 1494
$data[74,0] = "KovaaK's"
$data[74,1] = 32435
$data[75,0] = "Banana Shooter"
$data[75,1] = 968
$data[76,0] = "Don't Starve"
$data[76,1] = 3200
$data[77,0] = "CarX Drift Racing Online"
$data[77,1] = 870
$data[78,0] = "Teardown"
$data[78,1] = 6652
$data[79,0] = "Football Manager 2020"
$data[79,1] = 20132
$data[80,0] = "Barotrauma"
$data[80,1] = 54353
$data[81,0] = "Stranded: Alien Dawn"
$data[81,1] = 786
$data[82,0] = "The Elder Scrolls V: Skyrim"
$data[82,1] = 27731
$data[83,0] = "Songs of Syx"
$data[83,1] = 236
$data[84,0] = "House Flipper"
$data[84,1] = 30089
$data[85,0] = "Workers & Resources: Soviet Republic"
$data[85,1] = 9163
$data[86,0] = "Hydroneer"
$data[86,1] = 57
$data[87,0] = "Portal 2"
$data[87,1] = 948465
$data[88,0] = "Fisher Online"
$data[88,1] = 578
$data[89,0] = "PlateUp!"
$data[89,1] = 261
$data[90,0] = "Trove"
$data[90,1] = 2188
$data[91,0] = "Age of Mythology: Extended Edition"
$data[91,1] = 2223
$data[92,0] = "Wobbly Life"
$data[92,1] = 59
$data[93,0] = "Kerbal Space Program"
$data[93,1] = 105595
$data[94,0] = "Planet Coaster"
$data[94,1] = 408730
$data[95,0] = "Library Of Ruina"
$data[95,1] = 5920
$data[96,0] = "Scrap Mechanic"
$data[96,1] = 480190
$data[97,0] = "Stormworks: Build and Rescue"
$data[97,1] = 246550
$data[98,0] = "Age of Wonders 4"
$data[98,1] = 719
$data[99,0] = "Company of Heroes"
$data[99,1] = 3337

$ws.Range("B2:C101").Value = $data

